$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# Force column D cells to keep their original text representation (some
# values look numeric, e.g. "1.00", "497.30") by applying a temporary text
# number format before assigning the value, then restoring the default style
# so no stray formatting is left behind.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.248.95'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.277.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '497.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +2.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.151'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('E11').Value = '  +2.87%  '
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.682.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '54.198.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.273.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('E18').Value = '  +4.22%  '
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '304.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  +2.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.30'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '175.36'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.60%  '
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.95'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.51%  '
$ws.Range('E30').Value = '  +0.74%  '
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.927'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +10.47%  '
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.73'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.373'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '124.99'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('E42').Value = '  -1.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0492'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0896'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '239.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('E49').Value = '  +1.02%  '
$ws.Range('E50').Value = '  +0.31%  '
$ws.Range('E51').Value = '  +0.35%  '
